$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit re-orders the weekly price records (rows 2-17) onto a different
# set of dates/qualities/prices while columns A,B,C,E,F,G,H,N,O,Q,R stay put.
# Each new row's D/I/J/K/L/M/P values equal another row's original values
# (a permutation of the 16 data rows), so capture the "before" values first
# and then write them back out in the new order.

$cols = @("D","I","J","K","L","M","P")

# Snapshot current values for rows 2..17, columns D,I,J,K,L,M,P
# (Value2 is used for reads - this runtime's .Value getter does not behave
# as expected, but .Value2 round-trips both numbers and strings correctly.)
$snapshot = @{}
for ($r = 2; $r -le 17; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Destination row -> source row (source row's original values get copied in)
$mapping = @{
    2  = 13
    3  = 14
    4  = 15
    5  = 2
    6  = 3
    7  = 4
    8  = 5
    9  = 11
    10 = 16
    11 = 17
    12 = 10
    13 = 8
    14 = 9
    15 = 6
    16 = 7
    17 = 12
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
